$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B87: change from text "4" to a real number 4
$ws.Range("B87").Value = 4

# Add new row 88 with data
$ws.Range("A88").Value = "Ying Tang"

# B88 must be stored as TEXT "1" (not a number), even though it looks numeric.
# Assigning a plain numeric-looking string via .Value triggers Excel's
# quote-prefix (leading apostrophe) behavior, which stamps the cell with a
# new, distinct style. To keep the cell on the default style (like the rest
# of the row) we build the text via a TEXT() formula in a scratch cell, then
# copy/paste-special as values only, which yields a genuine text cell with
# no style change.
$ws.Range("Z1").Formula = '=TEXT(1,"0")'
$ws.Range("Z1").Copy()
$ws.Range("B88").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C88").Value = "trivial,of no interest at all,"
$ws.Range("D88").Value = "CRT"
$ws.Range("E88").Value = "THE"
$ws.Range("F88").Value = "91bd20d3-fd62-421d-99db-9f741ef9e1c2"
$ws.Range("G88").Value = "BkrsAzWAb_annotated.xlsx"
$ws.Range("H88").Value = "In contrast, the present result is trivial and of no interest at all, since it requires knowing a good parameter setting, which defeats a large part of the value of the method."
